$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1423.3334
$ws.Range("I98").Value = 1230.6666
$ws.Range("K98").Value = 1230.6666
$ws.Range("M98").Value = 267.3334

$ws.Range("H122").Value = 1423.3334
$ws.Range("I122").Value = 1230.6666
$ws.Range("K122").Value = 3691.9998
$ws.Range("M122").Value = -1241.9998

$ws.Range("H132").Value = 994
$ws.Range("I132").Value = 888.5185
$ws.Range("K132").Value = 2665.5555
$ws.Range("M132").Value = -135.5554999999999

$ws.Range("H141").Value = 4673886
$ws.Range("I141").Value = 7003274
$ws.Range("K141").Value = 21009822
$ws.Range("M141").Value = -21004642

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4966.289
$ws.Range("I32").Value = 4146.2197
$ws.Range("K32").Value = 4146.2197
$ws.Range("M32").Value = -3859.2197

$ws.Range("H46").Value = 10000
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()

$ws.Range("H61").Value = 6811.3335
$ws.Range("I61").Value = 7109.4
$ws.Range("J61").Value = 6066.1665
$ws.Range("K61").Value = 7109.4
$ws.Range("L61").Value = 6066.1665
$ws.Range("M61").Value = -6897.4
$ws.Range("N61").Value = -6490.1665

$ws.Range("H74").Value = 1358.4546
$ws.Range("I74").Value = 464.05884
$ws.Range("K74").Value = 464.05884
$ws.Range("M74").Value = 409.94116

$ws.Range("H77").Value = 1358.4546
$ws.Range("I77").Value = 464.05884
$ws.Range("K77").Value = 2320.2942
$ws.Range("M77").Value = 2047.7058

$ws.Range("H132").Value = 1149.2709
$ws.Range("I132").Value = 1045.1364
$ws.Range("K132").Value = 3135.4092
$ws.Range("M132").Value = -605.4092000000001

$ws.Range("H136").Value = 6811.3335
$ws.Range("I136").Value = 7109.4
$ws.Range("J136").Value = 6066.1665
$ws.Range("K136").Value = 21328.2
$ws.Range("L136").Value = 18198.4995
$ws.Range("M136").Value = -18778.2
$ws.Range("N136").Value = -23298.4995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 880.8333
$ws.Range("I22").Value = 821.75
$ws.Range("K22").Value = 821.75
$ws.Range("M22").Value = -648.75

$ws.Range("H99").Value = 1206.1875
$ws.Range("I99").Value = 1079.4
$ws.Range("K99").Value = 1079.4
$ws.Range("M99").Value = 418.5999999999999

$ws.Range("H134").Value = 2003.4231
$ws.Range("I134").Value = 1713.3158
$ws.Range("K134").Value = 5139.9474
$ws.Range("M134").Value = -2604.9474

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2192.9
$ws.Range("I31").Value = 1856.2858
$ws.Range("J31").Value = 2374.1538
$ws.Range("K31").Value = 1856.2858
$ws.Range("L31").Value = 2374.1538
$ws.Range("M31").Value = -1561.2858
$ws.Range("N31").Value = -2964.1538

$ws.Range("H34").Value = 2192.9
$ws.Range("I34").Value = 1856.2858
$ws.Range("J34").Value = 2374.1538
$ws.Range("K34").Value = 1856.2858
$ws.Range("L34").Value = 2374.1538
$ws.Range("M34").Value = -1654.2858
$ws.Range("N34").Value = -2778.1538

$ws.Range("H132").Value = 1609.6364
$ws.Range("I132").Value = 1119.9
$ws.Range("K132").Value = 3359.7
$ws.Range("M132").Value = -829.7000000000003

$ws.Range("H134").Value = 1477.303
$ws.Range("I134").Value = 1463.6562
$ws.Range("K134").Value = 4390.9686
$ws.Range("M134").Value = -1855.9686

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 18524.2
$ws.Range("J131").Value = 20262.414
$ws.Range("L131").Value = 60787.242
$ws.Range("N131").Value = -70867.242

$ws.Range("H140").Value = 2930.9546
$ws.Range("I140").Value = 1539.25
$ws.Range("J140").Value = 3240.2222
$ws.Range("K140").Value = 4617.75
$ws.Range("L140").Value = 9720.6666
$ws.Range("M140").Value = 562.25
$ws.Range("N140").Value = -20080.6666

$ws.Range("H141").Value = 4697.6665
$ws.Range("I141").Value = 5213.25
$ws.Range("J141").Value = 3666.5
$ws.Range("K141").Value = 15639.75
$ws.Range("L141").Value = 10999.5
$ws.Range("M141").Value = -10459.75
$ws.Range("N141").Value = -21359.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 12356.111
$ws.Range("I70").Value = 16467.5
$ws.Range("K70").Value = 16467.5
$ws.Range("M70").Value = -16197.5

$ws.Range("H73").Value = 12356.111
$ws.Range("I73").Value = 16467.5
$ws.Range("K73").Value = 16467.5
$ws.Range("M73").Value = -15531.5

$ws.Range("H102").Value = 1661.258
$ws.Range("I102").Value = 1348.625
$ws.Range("K102").Value = 1348.625
$ws.Range("M102").Value = 273.375

$ws.Range("H113").Value = 1495.6
$ws.Range("I113").Value = 1254
$ws.Range("J113").Value = 1556
$ws.Range("K113").Value = 1254
$ws.Range("L113").Value = 1556
$ws.Range("M113").Value = 916
$ws.Range("N113").Value = -5896

$ws.Range("H132").Value = 2748431.2
$ws.Range("I132").Value = 2959656.8
$ws.Range("K132").Value = 8878970.399999999
$ws.Range("M132").Value = -8876440.399999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1802.7858
$ws.Range("I46").Value = 1072.5
$ws.Range("K46").Value = 1072.5
$ws.Range("M46").Value = -884.5

$ws.Range("H61").Value = 2541.2856
$ws.Range("I61").Value = 2475.4285
$ws.Range("J61").Value = 2673
$ws.Range("K61").Value = 2475.4285
$ws.Range("L61").Value = 2673
$ws.Range("M61").Value = -2273.4285
$ws.Range("N61").Value = -3077

$ws.Range("H113").Value = 2541.2856
$ws.Range("I113").Value = 2475.4285
$ws.Range("J113").Value = 2673
$ws.Range("K113").Value = 2475.4285
$ws.Range("L113").Value = 2673
$ws.Range("M113").Value = -305.4285
$ws.Range("N113").Value = -7013

$ws.Range("H136").Value = 2621.2778
$ws.Range("I136").Value = 2523.3333
$ws.Range("K136").Value = 7569.999899999999
$ws.Range("M136").Value = -5019.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H122").Value = 27522.2
$ws.Range("I122").Value = 53513.934
$ws.Range("J122").Value = 1530.4667
$ws.Range("K122").Value = 160541.802
$ws.Range("L122").Value = 4591.4001
$ws.Range("M122").Value = -158091.802
$ws.Range("N122").Value = -9491.400099999999

$ws.Range("H132").Value = 1763.5294
$ws.Range("I132").Value = 1340.7916
$ws.Range("K132").Value = 4022.3748
$ws.Range("M132").Value = -1492.3748
